# "added ifoCAST full series evaluation"
# The error-table rows for horizons 6..15 (rows 2..11, sheet "first") are
# refreshed: each row's bias/MAE/RMSE/etc. figures (cols B:F) now reflect the
# next observation in the evaluation window (effectively the series shifted
# up by one row as a new period of actuals became available), a brand new
# figure is computed for the last row (row 11), and the "observations used"
# count in col G drops by one for every horizon.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (horizon 6)
$ws.Range("B2").Value = -0.0664853343312105
$ws.Range("C2").Value = 2.410518245823468
$ws.Range("D2").Value = 18.500884502483
$ws.Range("E2").Value = 4.301265453617459
$ws.Range("F2").Value = 4.401959587781231
$ws.Range("G2").Value = 22

# Row 3 (horizon 7)
$ws.Range("B3").Value = -0.5622320716953093
$ws.Range("C3").Value = 1.875594170520679
$ws.Range("D3").Value = 10.40438817800559
$ws.Range("E3").Value = 3.225583385684764
$ws.Range("F3").Value = 3.25464244416202
$ws.Range("G3").Value = 21

# Row 4 (horizon 8)
$ws.Range("B4").Value = -0.1812427017839789
$ws.Range("C4").Value = 1.441856625662851
$ws.Range("D4").Value = 8.636573577577872
$ws.Range("E4").Value = 2.938804787252442
$ws.Range("F4").Value = 3.009410619098892
$ws.Range("G4").Value = 20

# Row 5 (horizon 9)
$ws.Range("B5").Value = -0.2071371915338822
$ws.Range("C5").Value = 1.880800794843898
$ws.Range("D5").Value = 13.16300375377707
$ws.Range("E5").Value = 3.628085411587917
$ws.Range("F5").Value = 3.721423419543999
$ws.Range("G5").Value = 19

# Row 6 (horizon 10)
$ws.Range("B6").Value = -0.2618932902003571
$ws.Range("C6").Value = 1.772026561808276
$ws.Range("D6").Value = 9.805127963275318
$ws.Range("E6").Value = 3.131314095276186
$ws.Range("F6").Value = 3.21080636400065
$ws.Range("G6").Value = 18

# Row 7 (horizon 11)
$ws.Range("B7").Value = -0.2336496272831486
$ws.Range("C7").Value = 1.757688402756847
$ws.Range("D7").Value = 11.70993972551334
$ws.Range("E7").Value = 3.4219789195016
$ws.Range("F7").Value = 3.519063341396083
$ws.Range("G7").Value = 17

# Row 8 (horizon 12)
$ws.Range("B8").Value = -0.110650871936715
$ws.Range("C8").Value = 1.929024814611791
$ws.Range("D8").Value = 10.40150795858265
$ws.Range("E8").Value = 3.225136889898264
$ws.Range("F8").Value = 3.328946074360039
$ws.Range("G8").Value = 16

# Row 9 (horizon 13)
$ws.Range("B9").Value = -0.1706632301364332
$ws.Range("C9").Value = 2.15590519481267
$ws.Range("D9").Value = 15.75989337962914
$ws.Range("E9").Value = 3.969873219591419
$ws.Range("F9").Value = 4.105410294639327
$ws.Range("G9").Value = 15

# Row 10 (horizon 14)
$ws.Range("B10").Value = -0.1306111199879883
$ws.Range("C10").Value = 1.99991372347226
$ws.Range("D10").Value = 10.84813572458541
$ws.Range("E10").Value = 3.293650820075713
$ws.Range("F10").Value = 3.415294451084287
$ws.Range("G10").Value = 14

# Row 11 (horizon 15) - new data point for this horizon
$ws.Range("B11").Value = -0.1356576030215939
$ws.Range("C11").Value = 2.218579271193636
$ws.Range("D11").Value = 15.60620162443728
$ws.Range("E11").Value = 3.950468532267696
$ws.Range("F11").Value = 4.10935297329292
$ws.Range("G11").Value = 13
